$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 311.53845
$ws.Range("I12").Value = 206
$ws.Range("J12").Value = 549
$ws.Range("K12").Value = 206
$ws.Range("L12").Value = 549
$ws.Range("M12").Value = -36
$ws.Range("N12").Value = -889
$ws.Range("H18").Value = 774.75
$ws.Range("I18").Value = 774.75
$ws.Range("K18").Value = 774.75
$ws.Range("M18").Value = -490.75
$ws.Range("H19").Value = 2592
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 2592
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 2592
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -2942
$ws.Range("H32").Value = 4868.2144
$ws.Range("J32").Value = 5663.5
$ws.Range("L32").Value = 5663.5
$ws.Range("N32").Value = -6315.5
$ws.Range("H33").Value = 878
$ws.Range("I33").Value = 1080.5
$ws.Range("J33").Value = 704.4286
$ws.Range("K33").Value = 1080.5
$ws.Range("L33").Value = 704.4286
$ws.Range("M33").Value = -851.5
$ws.Range("N33").Value = -1162.4286
$ws.Range("H41").Value = 843.1667
$ws.Range("I41").Value = 596
$ws.Range("J41").Value = 1337.5
$ws.Range("K41").Value = 596
$ws.Range("L41").Value = 1337.5
$ws.Range("M41").Value = -156
$ws.Range("N41").Value = -2217.5
$ws.Range("H43").Value = 4187
$ws.Range("J43").Value = 4593.5
$ws.Range("L43").Value = 4593.5
$ws.Range("N43").Value = -4731.5
$ws.Range("H70").Value = 4077.88
$ws.Range("J70").Value = 4361.0454
$ws.Range("L70").Value = 13083.1362
$ws.Range("N70").Value = -13623.1362
$ws.Range("H73").Value = 4077.88
$ws.Range("J73").Value = 4361.0454
$ws.Range("L73").Value = 13083.1362
$ws.Range("N73").Value = -14955.1362
$ws.Range("H80").Value = 950.4
$ws.Range("I80").Value = 718.5
$ws.Range("J80").Value = 1105
$ws.Range("K80").Value = 2155.5
$ws.Range("L80").Value = 3315
$ws.Range("M80").Value = -1157.5
$ws.Range("N80").Value = -5311
$ws.Range("H83").Value = 950.4
$ws.Range("I83").Value = 718.5
$ws.Range("J83").Value = 1105
$ws.Range("K83").Value = 6466.5
$ws.Range("L83").Value = 9945
$ws.Range("M83").Value = -1474.5
$ws.Range("N83").Value = -19929
$ws.Range("H86").Value = 4885.143
$ws.Range("I86").Value = 4459.2
$ws.Range("K86").Value = 4459.2
$ws.Range("M86").Value = -3336.2
$ws.Range("H88").Value = 4151
$ws.Range("J88").Value = 4243
$ws.Range("L88").Value = 4243
$ws.Range("N88").Value = -5055
$ws.Range("H89").Value = 4885.143
$ws.Range("I89").Value = 4459.2
$ws.Range("K89").Value = 22296
$ws.Range("M89").Value = -16680
$ws.Range("H91").Value = 4151
$ws.Range("J91").Value = 4243
$ws.Range("L91").Value = 4243
$ws.Range("N91").Value = -7051
$ws.Range("H100").Value = 1568.5385
$ws.Range("I100").Value = 1113.875
$ws.Range("K100").Value = 1113.875
$ws.Range("M100").Value = -572.875
$ws.Range("H113").Value = 3654.8125
$ws.Range("I113").Value = 3697.7
$ws.Range("J113").Value = 3583.3333
$ws.Range("K113").Value = 3697.7
$ws.Range("L113").Value = 3583.3333
$ws.Range("M113").Value = -443.6999999999998
$ws.Range("N113").Value = -10091.3333
$ws.Range("H116").Value = 4889.8
$ws.Range("J116").Value = 4999.778
$ws.Range("L116").Value = 4999.778
$ws.Range("N116").Value = -11883.778
$ws.Range("H132").Value = 4458.794
$ws.Range("I132").Value = 1445.129
$ws.Range("J132").Value = 35600
$ws.Range("K132").Value = 4335.387
$ws.Range("L132").Value = 106800
$ws.Range("M132").Value = -1805.387
$ws.Range("N132").Value = -111860
$ws.Range("H135").Value = 990.0833
$ws.Range("I135").Value = 990.0833
$ws.Range("K135").Value = 8910.7497
$ws.Range("M135").Value = -6375.7497
$ws.Range("H137").Value = 2129.7856
$ws.Range("I137").Value = 1984.75
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 5954.25
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -3404.25
$ws.Range("N137").Value = -14100
$ws.Range("H138").Value = 2972.4243
$ws.Range("I138").Value = 3361.7778
$ws.Range("J138").Value = 2826.4167
$ws.Range("K138").Value = 10085.3334
$ws.Range("L138").Value = 8479.250100000001
$ws.Range("M138").Value = -4945.3334
$ws.Range("N138").Value = -18759.2501
$ws.Range("H141").Value = 2316.0833
$ws.Range("I141").Value = 2316.0833
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6948.249899999999
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -1768.249899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9844.923000000001
$ws.Range("I61").Value = 8712.809999999999
$ws.Range("J61").Value = 14599.8
$ws.Range("K61").Value = 8712.809999999999
$ws.Range("L61").Value = 14599.8
$ws.Range("M61").Value = -8500.809999999999
$ws.Range("N61").Value = -15023.8
$ws.Range("H74").Value = 5882.7
$ws.Range("I74").Value = 4883.923
$ws.Range("K74").Value = 4883.923
$ws.Range("M74").Value = -4009.923
$ws.Range("H77").Value = 5882.7
$ws.Range("I77").Value = 4883.923
$ws.Range("K77").Value = 24419.615
$ws.Range("M77").Value = -20051.615
$ws.Range("H136").Value = 9844.923000000001
$ws.Range("I136").Value = 8712.809999999999
$ws.Range("J136").Value = 14599.8
$ws.Range("K136").Value = 26138.43
$ws.Range("L136").Value = 43799.39999999999
$ws.Range("M136").Value = -23588.43
$ws.Range("N136").Value = -48899.39999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 20006
$ws.Range("I15").Value = 20006
$ws.Range("K15").Value = 20006
$ws.Range("M15").Value = -19779
$ws.Range("H20").Value = 4128.1514
$ws.Range("I20").Value = 3697.5557
$ws.Range("K20").Value = 3697.5557
$ws.Range("M20").Value = -3450.5557
$ws.Range("H22").Value = 651.2273
$ws.Range("I22").Value = 619.2941
$ws.Range("J22").Value = 759.8
$ws.Range("K22").Value = 619.2941
$ws.Range("L22").Value = 759.8
$ws.Range("M22").Value = -446.2941
$ws.Range("N22").Value = -1105.8
$ws.Range("H86").Value = 1916.7632
$ws.Range("I86").Value = 1726.258
$ws.Range("J86").Value = 2760.4285
$ws.Range("K86").Value = 1726.258
$ws.Range("L86").Value = 2760.4285
$ws.Range("M86").Value = -603.258
$ws.Range("N86").Value = -5006.4285
$ws.Range("H89").Value = 1916.7632
$ws.Range("I89").Value = 1726.258
$ws.Range("J89").Value = 2760.4285
$ws.Range("K89").Value = 8631.290000000001
$ws.Range("L89").Value = 13802.1425
$ws.Range("M89").Value = -3015.290000000001
$ws.Range("N89").Value = -25034.1425
$ws.Range("H102").Value = 4554.5
$ws.Range("I102").Value = 4554.5
$ws.Range("K102").Value = 4554.5
$ws.Range("M102").Value = -1309.5
$ws.Range("H107").Value = 2878.3333
$ws.Range("I107").Value = 2641.7144
$ws.Range("K107").Value = 2641.7144
$ws.Range("M107").Value = -721.7143999999998
$ws.Range("H134").Value = 5019.8
$ws.Range("I134").Value = 4832.25
$ws.Range("K134").Value = 14496.75
$ws.Range("M134").Value = -11961.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3406.1
$ws.Range("I16").Value = 1942.2
$ws.Range("K16").Value = 1942.2
$ws.Range("M16").Value = -1655.2
$ws.Range("H25").Value = 10000
$ws.Range("I25").Value = 10000
$ws.Range("K25").Value = 10000
$ws.Range("M25").Value = -9826
$ws.Range("H31").Value = 3604.8845
$ws.Range("I31").Value = 2684.6
$ws.Range("J31").Value = 4859.8184
$ws.Range("K31").Value = 2684.6
$ws.Range("L31").Value = 4859.8184
$ws.Range("M31").Value = -2389.6
$ws.Range("N31").Value = -5449.8184
$ws.Range("H34").Value = 3604.8845
$ws.Range("I34").Value = 2684.6
$ws.Range("J34").Value = 4859.8184
$ws.Range("K34").Value = 2684.6
$ws.Range("L34").Value = 4859.8184
$ws.Range("M34").Value = -2482.6
$ws.Range("N34").Value = -5263.8184
$ws.Range("H80").Value = 38327.797
$ws.Range("J80").Value = 38327.797
$ws.Range("L80").Value = 38327.797
$ws.Range("N80").Value = -40573.797
$ws.Range("H83").Value = 38327.797
$ws.Range("J83").Value = 38327.797
$ws.Range("L83").Value = 114983.391
$ws.Range("N83").Value = -126215.391
$ws.Range("H105").Value = 1688.6471
$ws.Range("I105").Value = 1282.6428
$ws.Range("J105").Value = 3583.3333
$ws.Range("K105").Value = 1282.6428
$ws.Range("L105").Value = 3583.3333
$ws.Range("M105").Value = 464.3571999999999
$ws.Range("N105").Value = -7077.3333
$ws.Range("H107").Value = 720.1667
$ws.Range("I107").Value = 571.8889
$ws.Range("K107").Value = 571.8889
$ws.Range("M107").Value = 1348.1111
$ws.Range("H113").Value = 3406.1
$ws.Range("I113").Value = 1942.2
$ws.Range("K113").Value = 1942.2
$ws.Range("M113").Value = 227.8
$ws.Range("H122").Value = 3914.125
$ws.Range("I122").Value = 4122.1113
$ws.Range("K122").Value = 12366.3339
$ws.Range("M122").Value = -9916.333899999998
$ws.Range("H132").Value = 4612.2
$ws.Range("I132").Value = 4824.154
$ws.Range("K132").Value = 14472.462
$ws.Range("M132").Value = -11942.462
$ws.Range("H141").Value = 37058
$ws.Range("I141").Value = 33822.5
$ws.Range("K141").Value = 33822.5
$ws.Range("M141").Value = -28642.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1560.7858
$ws.Range("J5").Value = 1857
$ws.Range("L5").Value = 5571
$ws.Range("N5").Value = -5795
$ws.Range("H107").Value = 381.92307
$ws.Range("J107").Value = 429.9091
$ws.Range("L107").Value = 1289.7273
$ws.Range("N107").Value = -5129.7273
$ws.Range("H113").Value = 1014.6667
$ws.Range("J113").Value = 1057.8
$ws.Range("L113").Value = 3173.4
$ws.Range("N113").Value = -7513.4
$ws.Range("H135").Value = 1560.7858
$ws.Range("J135").Value = 1857
$ws.Range("L135").Value = 16713
$ws.Range("N135").Value = -21783
$ws.Range("H138").Value = 2133.3333
$ws.Range("I138").Value = 2133.3333
$ws.Range("K138").Value = 6399.999899999999
$ws.Range("M138").Value = -1259.999899999999
$ws.Range("H140").Value = 2759.4546
$ws.Range("I140").Value = 2595.4443
$ws.Range("K140").Value = 7786.3329
$ws.Range("M140").Value = -2606.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 127742500
$ws.Range("J11").Value = 288000
$ws.Range("L11").Value = 288000
$ws.Range("N11").Value = -288278
$ws.Range("H113").Value = 403232.1
$ws.Range("I113").Value = 573640.9
$ws.Range("J113").Value = 5611.6665
$ws.Range("K113").Value = 573640.9
$ws.Range("L113").Value = 5611.6665
$ws.Range("M113").Value = -571470.9
$ws.Range("N113").Value = -9951.666499999999
$ws.Range("H122").Value = 5440.773
$ws.Range("J122").Value = 5353.7
$ws.Range("L122").Value = 16061.1
$ws.Range("N122").Value = -20961.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1710.4
$ws.Range("I16").Value = 1576.9412
$ws.Range("K16").Value = 1576.9412
$ws.Range("M16").Value = -1406.9412
$ws.Range("H122").Value = 3264.3684
$ws.Range("I122").Value = 2913.4
$ws.Range("J122").Value = 3654.3333
$ws.Range("K122").Value = 8740.200000000001
$ws.Range("L122").Value = 10962.9999
$ws.Range("M122").Value = -6290.200000000001
$ws.Range("N122").Value = -15862.9999
$ws.Range("H132").Value = 10566.625
$ws.Range("I132").Value = 11562.846
$ws.Range("J132").Value = 6249.6665
$ws.Range("K132").Value = 34688.538
$ws.Range("L132").Value = 18748.9995
$ws.Range("M132").Value = -32158.538
$ws.Range("N132").Value = -23808.9995
$ws.Range("H136").Value = 3468
$ws.Range("I136").Value = 3091.55
$ws.Range("K136").Value = 9274.650000000001
$ws.Range("M136").Value = -6724.650000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9996.666999999999
$ws.Range("J15").Value = 9995
$ws.Range("L15").Value = 9995
$ws.Range("N15").Value = -10571
$ws.Range("H100").Value = 779.25
$ws.Range("I100").Value = 1186.5
$ws.Range("J100").Value = 372
$ws.Range("K100").Value = 2373
$ws.Range("L100").Value = 744
$ws.Range("M100").Value = -1832
$ws.Range("N100").Value = -1826
$ws.Range("H107").Value = 1677
$ws.Range("I107").Value = 1216.3334
$ws.Range("K107").Value = 3649.0002
$ws.Range("M107").Value = -1729.0002
$ws.Range("H113").Value = 1090
$ws.Range("J113").Value = 450
$ws.Range("L113").Value = 1350
$ws.Range("N113").Value = -5690
$ws.Range("H132").Value = 3752.85
$ws.Range("I132").Value = 4043.0344
$ws.Range("J132").Value = 2987.818
$ws.Range("K132").Value = 12129.1032
$ws.Range("L132").Value = 8963.454000000002
$ws.Range("M132").Value = -9599.1032
$ws.Range("N132").Value = -14023.454
